$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.028.06"
$ws.Range("E2").Value = "  +3.28%  "

$ws.Range("D3").Value = "3.060.05"
$ws.Range("E3").Value = "  +2.29%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("D8").Value = "3.058.57"
$ws.Range("E8").Value = "  +2.36%  "

$ws.Range("E9").Value = "  +5.03%  "

$ws.Range("E10").Value = "  +6.19%  "

$ws.Range("E11").Value = "  -10.25%  "

$ws.Range("E12").Value = "  +9.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.77%  "

$ws.Range("D15").Value = "3.560.88"
$ws.Range("E15").Value = "  +2.24%  "

$ws.Range("D16").Value = "64.082.29"
$ws.Range("E16").Value = "  +3.18%  "

$ws.Range("D17").Value = "3.065.25"
$ws.Range("E17").Value = "  +2.19%  "

$ws.Range("E18").Value = "  +2.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.680"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +14.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.86%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("E27").Value = "  +3.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.92%  "

$ws.Range("E29").Value = "  +2.24%  "

$ws.Range("E30").Value = "  -0.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.84%  "

$ws.Range("E32").Value = "  +1.90%  "

$ws.Range("E33").Value = "  +4.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0408"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "444.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0811"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.28%  "

$ws.Range("D41").Value = "2.993.99"
$ws.Range("E41").Value = "  +2.26%  "

$ws.Range("E42").Value = "  +2.93%  "

$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.261"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.77%  "

$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.113"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.82%  "

$ws.Range("E49").Value = "  +5.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "118.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.05%  "
